# FUNCTIONALITY: Wrote two new test cases.
#
# Updates the Scheduling "Test Suite Statistics" sheet (Sheet1) with two
# newly-written, partially-automated test cases:
#   - Row 3 (SchedulingRequiredFields): status -> Finished, note about 1 case
#   - Row 4 (SchedulingUIFormat):       status -> Finished, note about 4 cases
# Totals (columns B/C) are bumped accordingly; the dependent SUM/ratio
# formulas in G4:G6 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: SchedulingRequiredFields ---------------------------------
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = "Finished"
$ws.Range("E3").Value = "Contains 1 partially automated test case."

# --- Row 4: SchedulingUIFormat ----------------------------------------
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 13
$ws.Range("D4").Value = "Finished"
$ws.Range("E4").Value = "Contains 4 partially automated test cases."

# --- Move the active selection to B5, matching the author's save state
$ws.Range("B5").Select()
